$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "Save" in H1 - copy G1's formatting (bold, border, centered) so it
# matches the rest of the header row, then set its text.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New "Save" data column: both data rows default to 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
